$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.9167996666666666
$ws.Range("H2").Value = 2.750399
$ws.Range("I2").Value = 0.2833456974325495
$ws.Range("J2").Value = 0.2833456974325495
$ws.Range("M2").Value = 6.799975
$ws.Range("N2").Value = 20.399925
$ws.Range("O2").Value = 0.2675773660825449
$ws.Range("P2").Value = 0.2675773660825449
$ws.Range("Q2").Value = 6.234214813341666
$ws.Range("R2").Value = 56.107933320075
$ws.Range("S2").Value = 0.0758168954098233
$ws.Range("T2").Value = 0.0758168954098233

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.9167996666666666
$ws.Range("H3").Value = 2.750399
$ws.Range("I3").Value = 0.2833456974325495
$ws.Range("J3").Value = 0.2833456974325495
$ws.Range("O3").Value = 0.283301948841328
$ws.Range("P3").Value = 0.283301948841328
$ws.Range("Q3").Value = 6.60057773933811
$ws.Range("R3").Value = 59.405199654043
$ws.Range("S3").Value = 0.08027238827844656
$ws.Range("T3").Value = 0.08027238827844656

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.9167996666666666
$ws.Range("H4").Value = 2.750399
$ws.Range("I4").Value = 0.2833456974325495
$ws.Range("J4").Value = 0.2833456974325495
$ws.Range("O4").Value = 0.4491206850761271
$ws.Range("P4").Value = 0.4491206850761271
$ws.Range("Q4").Value = 10.46394494748111
$ws.Range("R4").Value = 94.17550452732999
$ws.Range("S4").Value = 0.1272564137442797
$ws.Range("T4").Value = 0.1272564137442797

# Row 5
$ws.Range("I5").Value = 0.2271242616180895
$ws.Range("J5").Value = 0.2271242616180895
$ws.Range("M5").Value = 6.799975
$ws.Range("N5").Value = 20.399925
$ws.Range("O5").Value = 0.2675773660825449
$ws.Range("P5").Value = 0.2675773660825449
$ws.Range("Q5").Value = 4.997222294458333
$ws.Range("R5").Value = 44.975000650125
$ws.Range("S5").Value = 0.06077331169721123
$ws.Range("T5").Value = 0.06077331169721124

# Row 6
$ws.Range("I6").Value = 0.2271242616180895
$ws.Range("J6").Value = 0.2271242616180895
$ws.Range("O6").Value = 0.283301948841328
$ws.Range("P6").Value = 0.283301948841328
$ws.Range("S6").Value = 0.0643447459455524
$ws.Range("T6").Value = 0.0643447459455524

# Row 7
$ws.Range("I7").Value = 0.2271242616180895
$ws.Range("J7").Value = 0.2271242616180895
$ws.Range("O7").Value = 0.4491206850761271
$ws.Range("P7").Value = 0.4491206850761271
$ws.Range("S7").Value = 0.1020062039753259
$ws.Range("T7").Value = 0.1020062039753259

# Row 8
$ws.Range("G8").Value = 1.583934333333334
$ws.Range("H8").Value = 4.751803000000001
$ws.Range("I8").Value = 0.4895300409493609
$ws.Range("J8").Value = 0.4895300409493609
$ws.Range("M8").Value = 6.799975
$ws.Range("N8").Value = 20.399925
$ws.Range("O8").Value = 0.2675773660825449
$ws.Range("P8").Value = 0.2675773660825449
$ws.Range("Q8").Value = 10.77071386830833
$ws.Range("R8").Value = 96.93642481477501
$ws.Range("S8").Value = 0.1309871589755103
$ws.Range("T8").Value = 0.1309871589755103

# Row 9
$ws.Range("G9").Value = 1.583934333333334
$ws.Range("H9").Value = 4.751803000000001
$ws.Range("I9").Value = 0.4895300409493609
$ws.Range("J9").Value = 0.4895300409493609
$ws.Range("O9").Value = 0.283301948841328
$ws.Range("P9").Value = 0.283301948841328
$ws.Range("S9").Value = 0.1386848146173291
$ws.Range("T9").Value = 0.1386848146173291

# Row 10
$ws.Range("G10").Value = 1.583934333333334
$ws.Range("H10").Value = 4.751803000000001
$ws.Range("I10").Value = 0.4895300409493609
$ws.Range("J10").Value = 0.4895300409493609
$ws.Range("O10").Value = 0.4491206850761271
$ws.Range("P10").Value = 0.4491206850761271
$ws.Range("S10").Value = 0.2198580673565216
$ws.Range("T10").Value = 0.2198580673565216
